# Add blood and urine SOCs to all years
# (adds the two new "Subject of Collection" rows for Blood Sample / Urine
# Sample to the SSD sheet, following the same layout as the other SOC rows)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "SSD" sheet

# Set cell values in the same order the new shared strings were added in the
# original commit: ??blood, Blood Sample, ??urine, Urine Sample, then the two
# SOC URIs.
$ws1.Range("D10").Value = "??blood"
$ws1.Range("H10").Value = "Blood Sample"
$ws1.Range("D11").Value = "??urine"
$ws1.Range("H11").Value = "Urine Sample"
$ws1.Range("B10").Value = "nhanes-kb:SOC-NHANES-2009-2010-BLOOD"
$ws1.Range("B11").Value = "nhanes-kb:SOC-NHANES-2009-2010-URINE"

$ws1.Range("I10").Value = "nhanes-kb:STD-NHANES-2009-2010"
$ws1.Range("J10").Value = "nhanes-kb:SOC-NHANES-2009-2010-SUBJECTS"
$ws1.Range("I11").Value = "nhanes-kb:STD-NHANES-2009-2010"
$ws1.Range("J11").Value = "nhanes-kb:SOC-NHANES-2009-2010-SUBJECTS"

# Make the SSD sheet the active tab/sheet, and select rows 10:11 (the newly
# added rows), matching the saved selection/view state of the edited file.
$ws1.Activate()
$ws1.Range("A10:A11").EntireRow.Select()
